# "update brochure config excels" - refresh the AIBT / SISMIC brochure links:
# the Q2/courses-fees/new-region/non-CoE/SEAPAE rows are replaced by the
# current "AIBT Courses Fees 2021" and "AIBT Region2(SISMIC) Q4 Promotion"
# entries, and the sheet is trimmed back down to just those two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: AIBT course fees brochure
$ws.Range("A2").Value = "AIBT Courses Fees 2021.pdf"
$ws.Range("B2").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/sismic/aibt/AIBT_Courses_Fees_2021_VOL_2.2.pdf"

# Row 3: AIBT Region2 (SISMIC) Q4 promotion brochure
$ws.Range("A3").Value = "AIBT Region2(SISMIC) Q4 Promotion.pdf"
$ws.Range("B3").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/sismic/aibt/AIBTSISMIC_Q4_Brochure_1OCT-31DEC21_VOL1.1.pdf"

# The old extra brochure rows (4-6) are no longer published, drop them
$ws.Rows("4:6").Select()
$ws.Rows("4:6").Delete()

# Authoring machine's default workbook font, applied to the base/header styles
$wb.Styles.Item("Normal").Font.Name = "等线"
$ws.Range("A1:B1").Font.Bold = $true
